$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# --- Weekly crime-stat table updates (rows 15-31) ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 500
$ws.Range("N15").Value = 50
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = -5.882352941176
$ws.Range("L16").Value = 17.073170731707
$ws.Range("M16").Value = -20
$ws.Range("N16").Value = -88.292682926829
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 7.142857142857
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 66
$ws.Range("K17").Value = 56.060606060606
$ws.Range("L17").Value = 66.129032258064
$ws.Range("M17").Value = 194.285714285714
$ws.Range("N17").Value = 43.055555555555
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = -75
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 9.090909090909
$ws.Range("I18").Value = 93
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 55
$ws.Range("L18").Value = 29.166666666666
$ws.Range("M18").Value = 24
$ws.Range("N18").Value = -88.197969543147
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -6.25
$ws.Range("I19").Value = 215
$ws.Range("J19").Value = 273
$ws.Range("K19").Value = -21.245421245421
$ws.Range("L19").Value = -25.347222222222
$ws.Range("M19").Value = -0.921658986175
$ws.Range("N19").Value = -62.346760070052
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 94
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = -6
$ws.Range("L20").Value = 10.588235294117
$ws.Range("M20").Value = 70.90909090909
$ws.Range("N20").Value = -95.371738060068
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -8
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 2.564102564102
$ws.Range("I21").Value = 560
$ws.Range("J21").Value = 555
$ws.Range("K21").Value = 0.9009009009
$ws.Range("L21").Value = 1.633393829401
$ws.Range("M21").Value = 26.126126126126
$ws.Range("N21").Value = -85.563289507605
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 72.727272727272
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -27.5
$ws.Range("F24").Value = 136
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = 7.086614173228
$ws.Range("I24").Value = 1169
$ws.Range("J24").Value = 973
$ws.Range("K24").Value = 20.143884892086
$ws.Range("L24").Value = 30.178173719376
$ws.Range("M24").Value = 103.658536585366
$ws.Range("C25").Value = 19
$ws.Range("E25").Value = -26.923076923076
$ws.Range("F25").Value = 99
$ws.Range("G25").Value = 95
$ws.Range("H25").Value = 4.210526315789
$ws.Range("I25").Value = 910
$ws.Range("J25").Value = 702
$ws.Range("K25").Value = 29.629629629629
$ws.Range("L25").Value = 43.987341772151
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = -14.285714285714
$ws.Range("J26").Value = 167
$ws.Range("K26").Value = 0.598802395209
$ws.Range("L26").Value = 16.666666666666
$ws.Range("M26").Value = 25.373134328358
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = -9.090909090909
$ws.Range("L27").Value = 25
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -50
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = 20
$ws.Range("F31").Value = 1
$ws.Range("F31").NumberFormat = '#,##0'
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = 80
$ws.Range("L31").Value = -25
